$wb = $excel.ActiveWorkbook

# Rename the existing sheet and set up its selection (no active cell change,
# just select the whole data range and make sure it is no longer the
# selected/active tab once the new sheet is added below).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "by_prov"
$ws1.Range("A1:E11").Select()

# Add the new "all" sheet right after "by_prov".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "all"

# Populate the aggregated (all-province) OR summary table.
$ws2.Range("A1").Value = "group"
$ws2.Range("B1").Value = "OR"
$ws2.Range("C1").Value = "lower"
$ws2.Range("D1").Value = "higher"

$ws2.Range("A2").Value = "Anti-N Positivity"
$ws2.Range("B2").Value = 3.33
$ws2.Range("C2").Value = 2.78
$ws2.Range("D2").Value = 4.01

$ws2.Range("A3").Value = "Anti-S Positivity"
$ws2.Range("B3").Value = 1.27
$ws2.Range("C3").Value = 1.18
$ws2.Range("D3").Value = 1.36

# Match the column widths used for the equivalent columns on "by_prov"
# (20.42578125 / 8.42578125 "stored" character-width units; the closest the
# COM ColumnWidth setter's pixel grid can land on).
$ws2.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 7.666666666666667

# Leave the cursor parked just past the data, which is what Excel leaves
# behind after you've finished typing into the last populated cell.
$ws2.Range("D4").Select()

# "all" is the sheet the workbook should open on.
$ws2.Activate()
